# Fix the "Dia chi lien lac khi can bao tin" (address-for-notification) line:
# it was (incorrectly) bound to the ThiSinh_DienThoai (phone) merge field;
# rebind it to ThiSinh_DCNhanGiayBao (address to receive the admission letter).
#
# Commit message: "#update xuat phieu tong hop"

$d = $word.ActiveDocument

# Locate the target paragraph by its unique label text instead of a hard-coded
# paragraph index. There is a second, very similar-looking paragraph ("Dien
# thoai lien lac (de Truong bao tin): <<ThiSinh_DienThoai>>") that must stay
# untouched, so we match on the label that is unique to the paragraph we want.
$label = "li" + [char]0x00EA + "n l" + [char]0x1EA1 + "c khi c" + [char]0x1EA7 + "n b" + [char]0x00E1 + "o tin"

$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.Contains($label)) {
        $target = $para
        break
    }
}

$p = $target
$pStart = $p.Range.Start
$pEnd = $p.Range.End

# --- Step 1: rename the merge field placeholder text itself -----------------
$fieldRange = $d.Range($pStart, $pEnd)
$null = $fieldRange.Find.Execute("ThiSinh_DienThoai", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
$fieldNameLength = $fieldRange.End - $fieldRange.Start
$fieldRange.Text = "ThiSinh_DCNhanGiayBao"

# The replaced field text merged into a single run with no explicit font
# color; give it <w:color w:val="000000"/> like the other merge-field runs
# in this document (e.g. ThiSinh_HoKhauThuongTru).
$newFieldRange = $d.Range($fieldRange.Start, $fieldRange.Start + "ThiSinh_DCNhanGiayBao".Length)
$newFieldRange.Font.Color = 0

# --- Step 2: split the leading ": <<" run into ": " and "<<" ----------------
# After the rename above, the run holding ": <<" is untouched/merged as one
# run. Toggle Bold on/off on just the "<<" part: this forces the run engine
# to materialize a separate run boundary there, and turning Bold back off
# (its default) drops the <w:b/> marker again, leaving two plain runs with
# identical rPr - exactly matching the target XML.
$splitScope = $d.Range($pStart, $p.Range.End)
$null = $splitScope.Find.Execute("<<", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
$splitScope.Bold = 1
$splitScope.Bold = 0
